# Refresh the cryptos price/volume snapshot with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.515.76'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.602.58'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.16'
$ws.Range("E5").Value = '  +3.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.57'
$ws.Range("E6").Value = '  +1.53%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.50'
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("E11").Value = '  +1.39%  '
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.063.25'
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '59.419.02'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.86'
$ws.Range("E15").Value = '  +1.41%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.614.39'
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000133'
$ws.Range("E17").Value = '  +0.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '341.03'
$ws.Range("E18").Value = '  +0.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.37'
$ws.Range("E19").Value = '  +1.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.10'
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.35'
$ws.Range("E21").Value = '  -1.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.30'
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.165'
$ws.Range("E25").Value = '  -1.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  -0.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.21'
$ws.Range("E27").Value = '  +2.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0745'
$ws.Range("E28").Value = '  +2.71%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  +5.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.86'
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.82'
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.93'
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.840'
$ws.Range("E36").Value = '  +2.74%  '
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.825'
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '273.26'
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.601'
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.75'
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.54'
$ws.Range("E46").Value = '  +3.66%  '
$ws.Range("E47").Value = '  +1.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.940.22'
$ws.Range("E48").Value = '  -1.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.51'
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.34'
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.79'
$ws.Range("E51").Value = '  +1.40%  '

Write-Output "Updated cryptos sheet with latest values."
